$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("B12").Value = 60.5
$ws.Range("D12").Value = 99.83

# Row 78
$ws.Range("B78").Value = 67

# Row 85
$ws.Range("B85").Value = 167
$ws.Range("D85").Value = 233.8

# Row 131
$ws.Range("B131").Value = 6.5
$ws.Range("D131").Value = 119.4

# Row 182
$ws.Range("B182").Value = 88
$ws.Range("D182").Value = 264

# Row 213
$ws.Range("B213").Value = 41
$ws.Range("D213").Value = 196.8

# Row 244
$ws.Range("B244").Value = 60
$ws.Range("D244").Value = 250.8

# Row 288
$ws.Range("B288").Value = 28
$ws.Range("D288").Value = 200.98

# Row 368
$ws.Range("B368").Value = 23.5
$ws.Range("D368").Value = 184.24

# Row 457
$ws.Range("B457").Value = 17
$ws.Range("D457").Value = 127.5

# Row 579
$ws.Range("B579").Value = 6.5
$ws.Range("D579").Value = 12.75

# Row 581
$ws.Range("B581").Value = 117.5
$ws.Range("D581").Value = 333.45

# Row 605 (Grand Total)
$ws.Range("B605").Value = 36447.97
$ws.Range("D605").Value = 104231.78
